$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (SIQ ID = 1) updates
$ws.Range("E3").Value = "Floors number doesn't matter, because the motor will simulate the direction without the need for any counting "
$ws.Range("H3").Value = "23/1/2020"
$ws.Range("I3").Value = "TSH: 23/01/2020 : I agree with this we just need simulation for up and down movement"
$ws.Range("J3").Value = "Answered"

# Row 4 (SIQ ID = 2) updates
$ws.Range("H4").Value = "23/1/2020"
$ws.Range("I4").Value = "TSH: 23/01/2020 : Yes you should have a Keypad for input , it is great idea to have each key can navigate some characteres "
$ws.Range("J4").Value = "Answered"

# Row 5 (SIQ ID = 3) updates
$ws.Range("C5").Value = "`tThis requirement describe the reset functionality `nA reset for whole system is done when pressing on/off button for 2 sec"
$ws.Range("H5").Value = "23/1/2020"
$ws.Range("I5").Value = "TSH: 23/01/2020 : The reset here means that the system shall start as if it the first time , reset all user names and passwords and all start all over again "
$ws.Range("J5").Value = "Answered"

$ws.Rows.Item(5).RowHeight = 47.25

$ws.Columns.Item(4).ColumnWidth = 33.85546875
$ws.Columns.Item(9).ColumnWidth = 118.140625

$ws.Range("I9").Select()
